$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name and card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay as TEXT (matching the
# original inlineStr cell type), not be auto-converted to a number.
# Force text formatting before assignment, then restore the original
# cell formatting (NumberFormat="@" pulls in a "quote prefix" style)
# by re-pasting the formatting from an untouched sibling cell with the
# same base style (s="8").
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 31.08.2024"

# Row 6: transaction 1
$ws.Range("B6").Value = "01.09."
$ws.Range("C6").Value = "02.09."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-41716360"
$ws.Range("E6").Value = "53,04-"

# Row 7: transaction 2
$ws.Range("B7").Value = "04.09."
$ws.Range("C7").Value = "05.09."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU PJMWVI"
$ws.Range("E7").Value = "18,91-"

# Row 8: transaction 3
$ws.Range("B8").Value = "07.09."
$ws.Range("C8").Value = "08.09."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 91330493"
$ws.Range("E8").Value = "41,80-"

# Row 9: transaction 4
$ws.Range("B9").Value = "10.09."
$ws.Range("C9").Value = "11.09."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "25,46-"

# Row 12: closing balance date and amount
$ws.Range("D12").Value = "KONTOSTAND AM 15.09.2024"
$ws.Range("E12").Value = "139,21-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 25.09.2024"
